$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update flagged warnings: row 13 reverts to plain label, row 8 becomes flagged "No Error"
$ws.Range("O13").Value = "Flagged Warnings:"
$ws.Range("O8").Value = "Flagged Warnings: No Error"

# Row 2
$ws.Range("C2").Value = [double]"102.759488727824"
$ws.Range("D2").Value = [double]"0.01700636429314758"
$ws.Range("E2").Value = [double]"0.01203867778235281"
$ws.Range("F2").Value = [double]"1286.620747290223"
$ws.Range("G2").Value = [double]"0.008904149421854907"
$ws.Range("H2").Value = [double]"207.8816754167087"
$ws.Range("I2").Value = [double]"1286.620797292723"
$ws.Range("J2").Value = [double]"295.2300478297261"
$ws.Range("K2").Value = [double]"0.6427506593354835"
$ws.Range("L2").Value = [double]"2.221880731969836"
$ws.Range("M2").Value = [double]"0.1132357989236307"
$ws.Range("N2").Value = [double]"1.285501318670967"
$ws.Range("P2").Value = [double]"1389.380286020547"
$ws.Range("Q2").Value = [double]"0.008102214871292668"
$ws.Range("R2").Value = [double]"323.3896035716401"
$ws.Range("S2").Value = [double]"1389.380286020547"
$ws.Range("T2").Value = [double]"434.7075866717806"
$ws.Range("U2").Value = [double]"0.5937769230219699"
$ws.Range("V2").Value = [double]"0"
$ws.Range("W2").Value = [double]"2.92989190678327"
$ws.Range("X2").Value = [double]"0.1810746596166614"
$ws.Range("Y2").Value = [double]"1.18755384604394"
$ws.Range("AA2").Value = [double]"1266.367679775528"
$ws.Range("AB2").Value = [double]"34.25026094330259"
$ws.Range("AC2").Value = [double]"1.001057766623396"

# Row 3
$ws.Range("C3").Value = [double]"102.7355562076571"
$ws.Range("D3").Value = [double]"0.04182619887848488"
$ws.Range("E3").Value = [double]"0.03145488787956938"
$ws.Range("F3").Value = [double]"1286.70013705244"
$ws.Range("G3").Value = [double]"0.02848589665047804"
$ws.Range("H3").Value = [double]"45.6428136072575"
$ws.Range("I3").Value = [double]"1286.70013705244"
$ws.Range("J3").Value = [double]"69.63413326385562"
$ws.Range("K3").Value = [double]"0.7153551168605812"
$ws.Range("L3").Value = [double]"1.685106250405484"
$ws.Range("M3").Value = [double]"5.774065670038908E-09"
$ws.Range("N3").Value = [double]"1.430710233721162"
$ws.Range("P3").Value = [double]"1389.435693260097"
$ws.Range("Q3").Value = [double]"0.01334030222800683"
$ws.Range("R3").Value = [double]"84.35585673405524"
$ws.Range("S3").Value = [double]"1389.435693260097"
$ws.Range("T3").Value = [double]"111.8630356430008"
$ws.Range("U3").Value = [double]"0.5650093998242007"
$ws.Range("V3").Value = [double]"0"
$ws.Range("W3").Value = [double]"1.582587096471116"
$ws.Range("X3").Value = [double]"0.2751549582275368"
$ws.Range("Y3").Value = [double]"1.130018799648401"

# Row 4
$ws.Range("C4").Value = [double]"102.6855992771655"
$ws.Range("D4").Value = [double]"0.0395863000794863"
$ws.Range("E4").Value = [double]"0.02817174898886629"
$ws.Range("F4").Value = [double]"1286.704213272012"
$ws.Range("G4").Value = [double]"0.02204146764409983"
$ws.Range("H4").Value = [double]"62.52395601237291"
$ws.Range("I4").Value = [double]"1286.704213272012"
$ws.Range("J4").Value = [double]"85.7830367856584"
$ws.Range("K4").Value = [double]"0.6431807383925672"
$ws.Range("L4").Value = [double]"2.079424379651842"
$ws.Range("M4").Value = [double]"6.398916319039216E-08"
$ws.Range("N4").Value = [double]"1.286361476785134"
$ws.Range("P4").Value = [double]"1389.389812549177"
$ws.Range("Q4").Value = [double]"0.01754483243538646"
$ws.Range("R4").Value = [double]"96.52170514038116"
$ws.Range("S4").Value = [double]"1389.389812549177"
$ws.Range("T4").Value = [double]"127.8941887927423"
$ws.Range("U4").Value = [double]"0.6194449365778671"
$ws.Range("V4").Value = [double]"0"
$ws.Range("W4").Value = [double]"2.142561708832046"
$ws.Range("X4").Value = [double]"1.875747834834129E-09"
$ws.Range("Y4").Value = [double]"1.238889873155734"

# Row 5
$ws.Range("C5").Value = [double]"103.2068065483809"
$ws.Range("D5").Value = [double]"0.0472946170456274"
$ws.Range("E5").Value = [double]"0.03488891649976827"
$ws.Range("F5").Value = [double]"1285.207333773021"
$ws.Range("G5").Value = [double]"0.03067746121608884"
$ws.Range("H5").Value = [double]"40.96476818598457"
$ws.Range("I5").Value = [double]"1285.207333773021"
$ws.Range("J5").Value = [double]"113.3230864492769"
$ws.Range("K5").Value = [double]"0.9399007441921636"
$ws.Range("L5").Value = [double]"1.393693044088751"
$ws.Range("M5").Value = [double]"0.8645476682774037"
$ws.Range("N5").Value = [double]"1.879801488384327"
$ws.Range("P5").Value = [double]"1388.414140321402"
$ws.Range("Q5").Value = [double]"0.01661715582953856"
$ws.Range("R5").Value = [double]"70.14300020930993"
$ws.Range("S5").Value = [double]"1388.414140321402"
$ws.Range("T5").Value = [double]"163.1214269701559"
$ws.Range("U5").Value = [double]"0.9312448849758083"
$ws.Range("V5").Value = [double]"0"
$ws.Range("W5").Value = [double]"1.479035737403275"
$ws.Range("X5").Value = [double]"0.4534780741459566"
$ws.Range("Y5").Value = [double]"1.862489769951617"

# Row 6
$ws.Range("C6").Value = [double]"102.7400535638137"
$ws.Range("D6").Value = [double]"0.012279461555679"
$ws.Range("E6").Value = [double]"0.008683886320934935"
$ws.Range("F6").Value = [double]"1286.697275506042"
$ws.Range("G6").Value = [double]"0.006046742643592262"
$ws.Range("H6").Value = [double]"179.7541288309286"
$ws.Range("I6").Value = [double]"1286.697275506042"
$ws.Range("J6").Value = [double]"235.0844764676349"
$ws.Range("K6").Value = [double]"0.6143020373870589"
$ws.Range("L6").Value = [double]"3.138012966572981"
$ws.Range("M6").Value = [double]"6.527587780913446E-06"
$ws.Range("N6").Value = [double]"1.228604074774118"
$ws.Range("P6").Value = [double]"1389.437329069856"
$ws.Range("Q6").Value = [double]"0.006232718912086743"
$ws.Range("R6").Value = [double]"290.006191203727"
$ws.Range("S6").Value = [double]"1389.437329069856"
$ws.Range("T6").Value = [double]"379.3411626090142"
$ws.Range("U6").Value = [double]"0.5816713640334403"
$ws.Range("V6").Value = [double]"0"
$ws.Range("W6").Value = [double]"2.686882185178654"
$ws.Range("X6").Value = [double]"0.1620861707008304"
$ws.Range("Y6").Value = [double]"1.163342728066881"
$ws.Range("AA6").Value = [double]"1266.354250388924"
$ws.Range("AB6").Value = [double]"15.73254021978991"
$ws.Range("AC6").Value = [double]"0.4902690603456354"

# Row 7
$ws.Range("C7").Value = [double]"102.7341841700854"
$ws.Range("D7").Value = [double]"0.01489178893835922"
$ws.Range("E7").Value = [double]"0.01069378376823456"
$ws.Range("F7").Value = [double]"1286.689828647406"
$ws.Range("G7").Value = [double]"0.008763908579392018"
$ws.Range("H7").Value = [double]"244.7405949832576"
$ws.Range("I7").Value = [double]"1286.689828647406"
$ws.Range("J7").Value = [double]"311.3188451094823"
$ws.Range("K7").Value = [double]"0.5960543106132401"
$ws.Range("L7").Value = [double]"3.751446303086947"
$ws.Range("M7").Value = [double]"1.754782930074583E-09"
$ws.Range("N7").Value = [double]"1.19210862122648"
$ws.Range("P7").Value = [double]"1389.424012817491"
$ws.Range("Q7").Value = [double]"0.006127880358967198"
$ws.Range("R7").Value = [double]"385.1044144990192"
$ws.Range("S7").Value = [double]"1389.424012817491"
$ws.Range("T7").Value = [double]"491.3372303967248"
$ws.Range("U7").Value = [double]"0.5854443127621647"
$ws.Range("V7").Value = [double]"0"
$ws.Range("W7").Value = [double]"3.092547136975868"
$ws.Range("X7").Value = [double]"0.07013379742628134"
$ws.Range("Y7").Value = [double]"1.170888625524329"

# Row 8
$ws.Range("C8").Value = [double]"102.7607765319449"
$ws.Range("D8").Value = [double]"0.002978303392156865"
$ws.Range("E8").Value = [double]"0.002978303392156865"
$ws.Range("F8").Value = [double]"1286.675380917213"
$ws.Range("G8").Value = [double]"0"
$ws.Range("H8").Value = [double]"219.9336749782716"
$ws.Range("I8").Value = [double]"1286.675380917213"
$ws.Range("J8").Value = [double]"287.0962389048356"
$ws.Range("K8").Value = [double]"0.6033008331943694"
$ws.Range("L8").Value = [double]"3.065414054505084"
$ws.Range("M8").Value = [double]"0.04203610722108919"
$ws.Range("N8").Value = [double]"1.206601666388739"
$ws.Range("P8").Value = [double]"1389.436157449158"
$ws.Range("Q8").Value = [double]"0.002978303392156865"
$ws.Range("R8").Value = [double]"347.8887883033174"
$ws.Range("S8").Value = [double]"1389.436157449158"
$ws.Range("T8").Value = [double]"437.7398423560691"
$ws.Range("U8").Value = [double]"0.5849182516501219"
$ws.Range("V8").Value = [double]"0"
$ws.Range("W8").Value = [double]"2.908524365767789"
$ws.Range("X8").Value = [double]"0.03210826693627733"
$ws.Range("Y8").Value = [double]"1.169836503300244"
$ws.Range("AD8").Value = [double]"1410.546170618303"
$ws.Range("AE8").Value = [double]"35.329373001608"
$ws.Range("AF8").Value = [double]"0.4978467489414233"

# Row 9
$ws.Range("C9").Value = [double]"102.759339824592"
$ws.Range("D9").Value = [double]"0.01132346262347602"
$ws.Range("E9").Value = [double]"0.008765006006579866"
$ws.Range("F9").Value = [double]"1286.683837432701"
$ws.Range("G9").Value = [double]"0.008183132450267132"
$ws.Range("H9").Value = [double]"183.6164970473465"
$ws.Range("I9").Value = [double]"1286.683837432701"
$ws.Range("J9").Value = [double]"240.6091196850258"
$ws.Range("K9").Value = [double]"0.6155138442406235"
$ws.Range("L9").Value = [double]"3.071700865154769"
$ws.Range("M9").Value = [double]"2.520785401971537E-06"
$ws.Range("N9").Value = [double]"1.231027688481247"
$ws.Range("P9").Value = [double]"1389.443227259793"
$ws.Range("Q9").Value = [double]"0.003140330173208889"
$ws.Range("R9").Value = [double]"302.1507598415413"
$ws.Range("S9").Value = [double]"1389.443177257293"
$ws.Range("T9").Value = [double]"374.1738743289091"
$ws.Range("U9").Value = [double]"0.5816960636650494"
$ws.Range("V9").Value = [double]"0"
$ws.Range("W9").Value = [double]"3.196005289867535"
$ws.Range("X9").Value = [double]"1.28774279017918E-08"
$ws.Range("Y9").Value = [double]"1.163392127330099"
$ws.Range("AA9").Value = [double]"1266.373141335894"
$ws.Range("AB9").Value = [double]"19.20339722595636"
$ws.Range("AC9").Value = [double]"0.5751140545496863"
$ws.Range("AD9").Value = [double]"1410.627181908422"
$ws.Range("AE9").Value = [double]"42.16060873693058"
$ws.Range("AF9").Value = [double]"0.1992712256593459"

# Row 10
$ws.Range("C10").Value = [double]"102.7278206848084"
$ws.Range("D10").Value = [double]"0.01951959861319457"
$ws.Range("E10").Value = [double]"0.01435175972674143"
$ws.Range("F10").Value = [double]"1286.707651143"
$ws.Range("G10").Value = [double]"0.01254059435873648"
$ws.Range("H10").Value = [double]"156.0745342173684"
$ws.Range("I10").Value = [double]"1286.707651143"
$ws.Range("J10").Value = [double]"208.5880239143283"
$ws.Range("K10").Value = [double]"0.6245138943869486"
$ws.Range("L10").Value = [double]"3.075344287110709"
$ws.Range("M10").Value = [double]"2.533084852984757E-12"
$ws.Range("N10").Value = [double]"1.249027788773897"
$ws.Range("P10").Value = [double]"1389.435471827808"
$ws.Range("Q10").Value = [double]"0.00697900425445809"
$ws.Range("R10").Value = [double]"263.7608375111804"
$ws.Range("S10").Value = [double]"1389.435471827808"
$ws.Range("T10").Value = [double]"322.3103659786261"
$ws.Range("U10").Value = [double]"0.5726887961043112"
$ws.Range("V10").Value = [double]"0"
$ws.Range("W10").Value = [double]"2.942496238938431"
$ws.Range("X10").Value = [double]"6.568381893945485E-08"
$ws.Range("Y10").Value = [double]"1.145377592208622"

# Row 11
$ws.Range("C11").Value = [double]"102.7296626981142"
$ws.Range("D11").Value = [double]"0.01891010925144206"
$ws.Range("E11").Value = [double]"0.01559057432361063"
$ws.Range("F11").Value = [double]"1286.714117859994"
$ws.Range("G11").Value = [double]"0.01512382397141193"
$ws.Range("H11").Value = [double]"130.8125811186033"
$ws.Range("I11").Value = [double]"1286.714117859994"
$ws.Range("J11").Value = [double]"173.7138312181267"
$ws.Range("K11").Value = [double]"0.6199414396133358"
$ws.Range("L11").Value = [double]"2.919465740733076"
$ws.Range("M11").Value = [double]"6.42167985454023E-11"
$ws.Range("N11").Value = [double]"1.239882879226672"
$ws.Range("P11").Value = [double]"1389.443780558108"
$ws.Range("Q11").Value = [double]"0.003786285280030128"
$ws.Range("R11").Value = [double]"235.099830632959"
$ws.Range("S11").Value = [double]"1389.443780558108"
$ws.Range("T11").Value = [double]"289.8211297001717"
$ws.Range("U11").Value = [double]"0.5165213034276208"
$ws.Range("V11").Value = [double]"0"
$ws.Range("W11").Value = [double]"2.608060848345964"
$ws.Range("X11").Value = [double]"0.3350006185512883"
$ws.Range("Y11").Value = [double]"1.033042606855242"
$ws.Range("AD11").Value = [double]"1410.669144030468"
$ws.Range("AE11").Value = [double]"23.11247378456391"
$ws.Range("AF11").Value = [double]"0.539675702924542"

# Row 12
$ws.Range("C12").Value = [double]"102.7583746593305"
$ws.Range("D12").Value = [double]"0.01958960327781861"
$ws.Range("E12").Value = [double]"0.01435943676881631"
$ws.Range("F12").Value = [double]"1286.687722802019"
$ws.Range("G12").Value = [double]"0.01247035261267834"
$ws.Range("H12").Value = [double]"126.1765636384594"
$ws.Range("I12").Value = [double]"1286.687722802019"
$ws.Range("J12").Value = [double]"160.3687132009049"
$ws.Range("K12").Value = [double]"0.5931448825945722"
$ws.Range("L12").Value = [double]"2.640592896090385"
$ws.Range("M12").Value = [double]"3.702593787124897E-14"
$ws.Range("N12").Value = [double]"1.186289765189144"
$ws.Range("P12").Value = [double]"1389.44614746385"
$ws.Range("Q12").Value = [double]"0.007119250665140266"
$ws.Range("R12").Value = [double]"211.8219846285743"
$ws.Range("S12").Value = [double]"1389.44609746135"
$ws.Range("T12").Value = [double]"265.3448276973129"
$ws.Range("U12").Value = [double]"0.5692989110669842"
$ws.Range("V12").Value = [double]"0"
$ws.Range("W12").Value = [double]"2.331539725587954"
$ws.Range("X12").Value = [double]"0.1008219568759909"
$ws.Range("Y12").Value = [double]"1.138597822133968"
$ws.Range("AD12").Value = [double]"1409.463237821099"
$ws.Range("AE12").Value = [double]"268.9485272698817"
$ws.Range("AF12").Value = [double]"0.2108563242014863"

# Row 13
$ws.Range("C13").Value = [double]"102.7274980038515"
$ws.Range("D13").Value = [double]"0.01935201509201873"
$ws.Range("E13").Value = [double]"0.01414613921773145"
$ws.Range("F13").Value = [double]"1286.718261236902"
$ws.Range("G13").Value = [double]"0.01221204882195297"
$ws.Range("H13").Value = [double]"118.3712069097287"
$ws.Range("I13").Value = [double]"1286.718261236902"
$ws.Range("J13").Value = [double]"152.9912346300638"
$ws.Range("K13").Value = [double]"0.6027299408373412"
$ws.Range("L13").Value = [double]"2.535498106193613"
$ws.Range("M13").Value = [double]"1.589177123229035E-10"
$ws.Range("N13").Value = [double]"1.205459881674682"
$ws.Range("P13").Value = [double]"1389.445759240753"
$ws.Range("Q13").Value = [double]"0.00713996627006576"
$ws.Range("R13").Value = [double]"194.7497236108655"
$ws.Range("S13").Value = [double]"1389.445759240753"
$ws.Range("T13").Value = [double]"240.7648242427668"
$ws.Range("U13").Value = [double]"0.5678424365943296"
$ws.Range("V13").Value = [double]"0"
$ws.Range("W13").Value = [double]"2.139305002554165"
$ws.Range("X13").Value = [double]"0.06464257500490034"
$ws.Range("Y13").Value = [double]"1.135684873188659"

# Clear cells that no longer appear in the final state
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AA4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("AD4").ClearContents()
$ws.Range("AE4").ClearContents()
$ws.Range("AF4").ClearContents()
$ws.Range("AD6").ClearContents()
$ws.Range("AE6").ClearContents()
$ws.Range("AF6").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AA10").ClearContents()
$ws.Range("AB10").ClearContents()
$ws.Range("AC10").ClearContents()
$ws.Range("AD10").ClearContents()
$ws.Range("AE10").ClearContents()
$ws.Range("AF10").ClearContents()
$ws.Range("AA11").ClearContents()
$ws.Range("AB11").ClearContents()
$ws.Range("AC11").ClearContents()
$ws.Range("AA12").ClearContents()
$ws.Range("AB12").ClearContents()
$ws.Range("AC12").ClearContents()
$ws.Range("AA13").ClearContents()
$ws.Range("AB13").ClearContents()
$ws.Range("AC13").ClearContents()
$ws.Range("AD13").ClearContents()
$ws.Range("AE13").ClearContents()
$ws.Range("AF13").ClearContents()

Write-Host "applied edits"
